$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 64476000000.0
$ws.Range("G2").Value = 54855000000.0

$ws.Range("B4").Value = 2570000000.0

$ws.Range("B6").Value = 57922000000.0
$ws.Range("G6").Value = 35323000000.0

$ws.Range("B7").Value = 6463000000.0
$ws.Range("G7").Value = 86000000.0

$ws.Range("B9").Value = 2255000000.0
$ws.Range("G9").Value = 2673000000.0

$ws.Range("B12").Value = 1562000000.0
$ws.Range("G12").Value = 1363000000.0

$ws.Range("B14").Value = 1040000000.0

$ws.Range("B17").Value = 9793000000.0

$ws.Range("B19").Value = 10574000000.0

$ws.Range("B20").Value = 418000000.0

$ws.Range("G21").Value = 7745000000.0

$ws.Range("B24").Value = 51160000000.0

$ws.Range("B31").Value = 2841000000.0

$ws.Range("B34").Value = -52862000000.0

$ws.Range("B35").Value = 11614000000.0
